$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") for data rows 2..56: 46070 -> 46072
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 46072
}

# 2) Rows 43 and 44 swap their A (Beteckning), B (Datum) and G (Area (ha)) values.
#    (columns D/E/H..R are identical between the two rows, so only A/B/G visibly move)
$a43 = $ws.Cells.Item(43, 1).Value2
$b43 = $ws.Cells.Item(43, 2).Value2
$g43 = $ws.Cells.Item(43, 7).Value2

$a44 = $ws.Cells.Item(44, 1).Value2
$b44 = $ws.Cells.Item(44, 2).Value2
$g44 = $ws.Cells.Item(44, 7).Value2

$ws.Cells.Item(43, 1).Value2 = $a44
$ws.Cells.Item(43, 2).Value2 = $b44
$ws.Cells.Item(43, 7).Value2 = $g44

$ws.Cells.Item(44, 1).Value2 = $a43
$ws.Cells.Item(44, 2).Value2 = $b43
$ws.Cells.Item(44, 7).Value2 = $g43
